# Clear out the contents of row 54 (A54:J54) on the active sheet, as if the
# user selected the row and pressed Delete, leaving only D54's date
# formatting (style) behind with no value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Range("A54:J54")
$row.Select()
$row.ClearContents()
